$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "purpose" column (E) for rows 2-25 from "S.GISH" to "fullRNASEQ"
$ws.Range("E2:E25").Value = "fullRNASEQ"

# Update the sheet view: scroll position + selection
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollRow = 17
$ws.Range("E24:E25").Select()

# Enable iterative calculation delta tuning
$excel.MaxChange = 0.0001
